$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the "MuSCs" sending-cluster row; rows below shift up (old row 6 -> new row 5).
$ws.Rows.Item(5).Delete()

# Refresh the TPM-derived statistics for the remaining data rows (2-5) with the new values.
$newValues = @{
    2 = @(2, 0.6666666666666666, 0.7644876666666667, 2.293463, 0.1573643627333495, 0.1573643627333495, 1, 0.3333333333333333, 0.08314566666666666, 0.249437, 1, 1, 0.06356383670344444, 0.572074530331, 0.1573643627333495, 0.1573643627333495)
    3 = @(2, 0.6666666666666666, 0.128779, 0.386337, 0.0265082435623832, 0.0265082435623832, 1, 0.3333333333333333, 0.08314566666666666, 0.249437, 1, 1, 0.01070741580766667, 0.09636674226900001, 0.0265082435623832, 0.0265082435623832)
    4 = @(3, 1, 1.975728666666667, 5.927186, 0.406689729763258, 0.406689729763258, 1, 0.3333333333333333, 0.08314566666666666, 0.249437, 1, 1, 0.1642732771424444, 1.478459494282, 0.406689729763258, 0.406689729763258)
    5 = @(3, 1, 1.989078333333334, 5.967235000000001, 0.4094376639410093, 0.4094376639410093, 1, 0.3333333333333333, 0.08314566666666666, 0.249437, 1, 1, 0.1653832440772222, 1.488449196695, 0.4094376639410093, 0.4094376639410093)
}

foreach ($r in $newValues.Keys) {
    $vals = $newValues[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 5 + $i
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}
